# "Set write direction through control panel" -
# Insert a new checkbox-style output-control row ("output_horizontally")
# into the "constants" sheet, just below "output_by_scenario" (row 81),
# pushing the remaining control-panel rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("constants")

# Insert a new row at 82; Excel/COM convention inherits formatting from
# the row above (row 81), matching the other checkbox rows in this block.
$ws.Rows.Item(82).Insert()

# Populate the new control-panel row: label + an unchecked boolean toggle.
$ws.Range("A82").Value = "output_horizontally"
$ws.Range("B82").Value = $false

# Match the author's final selection state after the edit.
$ws.Range("C81").Select() | Out-Null

Write-Host "Inserted 'output_horizontally' row at constants!A82"
